$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 141, pushing the existing rows 141-163
# (and everything after) down to 143-165. This mirrors the weekly update
# where a new date's data (serial 44505) was added at the top of the
# historical block.
$ws.Rows("141:142").Insert()

# Row 141: Packham's Triumph, Primera
$ws.Range("A141").Value2 = 4
$ws.Range("B141").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C141").Value2 = "Los Lagos"
$ws.Range("D141").Value2 = 44505
$ws.Range("E141").Value2 = 10
$ws.Range("F141").Value2 = "Fruta"
$ws.Range("G141").Value2 = 100104
$ws.Range("H141").Value2 = "Frutos de pepita"
$ws.Range("I141").Value2 = 100104005
$ws.Range("J141").Value2 = "Pera"
$ws.Range("K141").Value2 = "Packham's Triumph"
$ws.Range("L141").Value2 = "Primera"
$ws.Range("M141").Value2 = 500
$ws.Range("N141").Value2 = 15000
$ws.Range("O141").Value2 = 16000
$ws.Range("P141").Value2 = 15500
$ws.Range("Q141").Value2 = "`$/caja 15 kilos empedrada"
$ws.Range("R141").Value2 = "Región de O'Higgins"
$ws.Range("S141").Value2 = 1033
$ws.Range("T141").Value2 = 15

# Row 142: Packham's Triumph, Segunda
$ws.Range("A142").Value2 = 4
$ws.Range("B142").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C142").Value2 = "Los Lagos"
$ws.Range("D142").Value2 = 44505
$ws.Range("E142").Value2 = 10
$ws.Range("F142").Value2 = "Fruta"
$ws.Range("G142").Value2 = 100104
$ws.Range("H142").Value2 = "Frutos de pepita"
$ws.Range("I142").Value2 = 100104005
$ws.Range("J142").Value2 = "Pera"
$ws.Range("K142").Value2 = "Packham's Triumph"
$ws.Range("L142").Value2 = "Segunda"
$ws.Range("M142").Value2 = 200
$ws.Range("N142").Value2 = 14000
$ws.Range("O142").Value2 = 14000
$ws.Range("P142").Value2 = 14000
$ws.Range("Q142").Value2 = "`$/caja 15 kilos empedrada"
$ws.Range("R142").Value2 = "Región de O'Higgins"
$ws.Range("S142").Value2 = 933
$ws.Range("T142").Value2 = 15
